$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (which currently holds "Tipo").
# This shifts the existing "Tipo" column (and its data) from D to E.
$ws.Range("D1").EntireColumn.Insert()

# Set header for the newly inserted column D, matching the style of the
# other header cells (bold header with border), same as B1/C1/(old D1).
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "MAE"

# Fill in the MAE value for row 2
$ws.Range("D2").Value = 0.4321333824756282

# Update the slightly-changed MSE value in B2
$ws.Range("B2").Value = 0.3955935532374558
